# Swap the contents of columns B:AD between pairs of rows.
# Column A (the running id/index) is intentionally left untouched on both rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(18, 19),
    @(28, 29),
    @(52, 53),
    @(108, 109),
    @(125, 126),
    @(130, 131),
    @(133, 134),
    @(164, 165),
    @(197, 198),
    @(203, 204),
    @(210, 211),
    @(218, 219),
    @(226, 227),
    @(229, 230)
)

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    $rangeA = $ws.Range("B$rowA`:AD$rowA")
    $rangeB = $ws.Range("B$rowB`:AD$rowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}
